# Trade #185 closed at 2026-02-17 22:04:44 - unknown UNKNOWN +0.000%
#
# Updates the live trading results workbook:
#  - Summary sheet roll-up numbers
#  - Strategy Status row for MarketMaking
#  - Closes trade #213 (row 214 on "All Trades", row 181 on "MarketMaking")
#    with an early_exit reason
#  - Appends the two newest open trades (#246 volatility_scorer,
#    #247 MarketMaking) to "All Trades" and to their respective
#    per-strategy sheets

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value2 = 1400        # Current Capital
$summary.Range("B4").Value2 = -0.22       # Total P&L $
$summary.Range("B6").Value2 = 213         # Total Trades
$summary.Range("B8").Value2 = 93          # Losing Trades
$summary.Range("B9").Value2 = 38.5        # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value2 = 100          # Capital
$status.Range("D5").Value2 = 180          # Trades
$status.Range("E5").Value2 = -0.33        # P&L $
$status.Range("F5").Value2 = -0           # P&L %
$status.Range("G5").Value2 = 37.78        # Win Rate %

# ---------------------------------------------------------------------
# All Trades sheet - close trade #213 (row 214) + append new trades
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

$allTrades.Cells.Item(214, 7).Value2 = 0.01        # G214 Exit Price
$allTrades.Cells.Item(214, 8).Value2 = "CLOSED"    # H214 Status
$allTrades.Cells.Item(214, 9).Value2 = -75         # I214 P&L %
$allTrades.Cells.Item(214, 10).Value2 = -0.03      # J214 P&L $
$allTrades.Cells.Item(214, 11).Value2 = 100        # K214 Capital After
$allTrades.Cells.Item(214, 12).Value2 = "early_exit" # L214 Exit Reason
$allTrades.Cells.Item(214, 13).Value2 = 0.2        # M214 Duration (min)

# New row 247: trade #246 (volatility_scorer, still OPEN)
$allTrades.Cells.Item(247, 1).Value2 = 246
$allTrades.Cells.Item(247, 2).Value2 = "'2026-02-17"
$allTrades.Cells.Item(247, 3).Value2 = "'22:04:36"
$allTrades.Cells.Item(247, 4).Value2 = "volatility_scorer"
$allTrades.Cells.Item(247, 5).Value2 = "NEUTRAL"
$allTrades.Cells.Item(247, 6).Value2 = 0.04
$allTrades.Cells.Item(247, 7).Value2 = "'"          # G247 blank (exit price not set yet)
$allTrades.Cells.Item(247, 8).Value2 = "OPEN"
$allTrades.Cells.Item(247, 9).Value2 = 0
$allTrades.Cells.Item(247, 10).Value2 = 0
$allTrades.Cells.Item(247, 11).Value2 = 100
$allTrades.Cells.Item(247, 12).Value2 = "'"         # L247 blank (exit reason not set yet)
$allTrades.Cells.Item(247, 13).Value2 = 0
$allTrades.Cells.Item(247, 14).Value2 = 0
$allTrades.Cells.Item(247, 15).Value2 = 0
$allTrades.Cells.Item(247, 16).Value2 = 0.85
$allTrades.Cells.Item(247, 17).Value2 = "Low vol market (score: inf) - ideal for market making"

# New row 248: trade #247 (MarketMaking, still OPEN)
$allTrades.Cells.Item(248, 1).Value2 = 247
$allTrades.Cells.Item(248, 2).Value2 = "'2026-02-17"
$allTrades.Cells.Item(248, 3).Value2 = "'22:04:38"
$allTrades.Cells.Item(248, 4).Value2 = "MarketMaking"
$allTrades.Cells.Item(248, 5).Value2 = "DOWN"
$allTrades.Cells.Item(248, 6).Value2 = 0.04
$allTrades.Cells.Item(248, 7).Value2 = "'"          # G248 blank (exit price not set yet)
$allTrades.Cells.Item(248, 8).Value2 = "OPEN"
$allTrades.Cells.Item(248, 9).Value2 = 0
$allTrades.Cells.Item(248, 10).Value2 = 0
$allTrades.Cells.Item(248, 11).Value2 = 100.0269683756113
$allTrades.Cells.Item(248, 12).Value2 = "'"         # L248 blank (exit reason not set yet)
$allTrades.Cells.Item(248, 13).Value2 = 0
$allTrades.Cells.Item(248, 14).Value2 = 0
$allTrades.Cells.Item(248, 15).Value2 = 0
$allTrades.Cells.Item(248, 16).Value2 = 0.6
$allTrades.Cells.Item(248, 17).Value2 = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------
# volatility_scorer sheet - append new row 5 (trade #246)
# ---------------------------------------------------------------------
$volScorer = $wb.Worksheets.Item("volatility_scorer")

$volScorer.Cells.Item(5, 1).Value2 = 246
$volScorer.Cells.Item(5, 2).Value2 = "'2026-02-17"
$volScorer.Cells.Item(5, 3).Value2 = "'22:04:36"
$volScorer.Cells.Item(5, 4).Value2 = "volatility_scorer"
$volScorer.Cells.Item(5, 5).Value2 = "NEUTRAL"
$volScorer.Cells.Item(5, 6).Value2 = 0.04
$volScorer.Cells.Item(5, 7).Value2 = "'"            # G5 blank (exit price not set yet)
$volScorer.Cells.Item(5, 8).Value2 = "OPEN"
$volScorer.Cells.Item(5, 9).Value2 = 0
$volScorer.Cells.Item(5, 10).Value2 = 0
$volScorer.Cells.Item(5, 11).Value2 = 100
$volScorer.Cells.Item(5, 12).Value2 = 0
$volScorer.Cells.Item(5, 13).Value2 = 0
$volScorer.Cells.Item(5, 14).Value2 = 0.85
$volScorer.Cells.Item(5, 15).Value2 = "Low vol market (score: inf) - ideal for market making"
$volScorer.Cells.Item(5, 16).Value2 = "'"           # P5 blank (exit reason not set yet)
$volScorer.Cells.Item(5, 17).Value2 = 0

# ---------------------------------------------------------------------
# MarketMaking sheet - close trade #213 (row 181) + append new row 211
# ---------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")

$marketMaking.Cells.Item(181, 7).Value2 = 0.01         # G181 Exit Price
$marketMaking.Cells.Item(181, 8).Value2 = "CLOSED"     # H181 Status
$marketMaking.Cells.Item(181, 9).Value2 = -75          # I181 P&L %
$marketMaking.Cells.Item(181, 10).Value2 = -0.03       # J181 P&L $
$marketMaking.Cells.Item(181, 11).Value2 = 100         # K181 Capital After
$marketMaking.Cells.Item(181, 16).Value2 = "early_exit"  # P181 Exit Reason
$marketMaking.Cells.Item(181, 17).Value2 = 0.2         # Q181 Duration (min)

# New row 211: trade #247 (MarketMaking, still OPEN)
$marketMaking.Cells.Item(211, 1).Value2 = 247
$marketMaking.Cells.Item(211, 2).Value2 = "'2026-02-17"
$marketMaking.Cells.Item(211, 3).Value2 = "'22:04:38"
$marketMaking.Cells.Item(211, 4).Value2 = "MarketMaking"
$marketMaking.Cells.Item(211, 5).Value2 = "DOWN"
$marketMaking.Cells.Item(211, 6).Value2 = 0.04
$marketMaking.Cells.Item(211, 7).Value2 = "'"       # G211 blank (exit price not set yet)
$marketMaking.Cells.Item(211, 8).Value2 = "OPEN"
$marketMaking.Cells.Item(211, 9).Value2 = 0
$marketMaking.Cells.Item(211, 10).Value2 = 0
$marketMaking.Cells.Item(211, 11).Value2 = 100.0269683756113
$marketMaking.Cells.Item(211, 12).Value2 = 0
$marketMaking.Cells.Item(211, 13).Value2 = 0
$marketMaking.Cells.Item(211, 14).Value2 = 0.6
$marketMaking.Cells.Item(211, 15).Value2 = "Normal spread capture: 19600 bps"
$marketMaking.Cells.Item(211, 16).Value2 = "'"      # P211 blank (exit reason not set yet)
$marketMaking.Cells.Item(211, 17).Value2 = 0
